$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sigmas")
$ws.Activate()

# --- Cell value updates ---

# Row 2
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("K3").ClearContents()

# Row 5
$ws.Range("B5").Value = 1

# Row 10
$ws.Range("D10").Value = 2
$ws.Range("F10").Value = 2

# Row 11
$ws.Range("D11").Value = 1

# Row 15
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 2

# Row 25
$ws.Range("D25").Value = 2
$ws.Range("F25").Value = 2

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = 2

# --- View state updates ---
# Move the active selection to E4 (frozen-pane scroll position itself is
# session/host view state that this runtime does not persist independent of
# the freeze split, so only the selected cell is updated here).
$ws.Range("E4").Select()
